# Applies the "Update front page and Italy" edit to destroyers.docx
# 1. Italy intro paragraph: "used" -> "uses"
# 2. Audace Class paragraph: "Four OTO-" -> "A single 127mm and four OTO-"
# 3. De la Penne Class paragraph: "are were designed to be improved" -> "were designed as improved"
# 4. Same paragraph: armament sentence rewritten (three 76mm guns vice four, no more Audace mention)
# 5. Move the _GoBack bookmark from the trailing empty paragraph to right after "vice four"

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING - text not found: $findText"
    }
    return $ok
}

# Change 1
$find1 = "four destroyers and used them"
$repl1 = "four destroyers and uses them"
Replace-Text $find1 $repl1 | Out-Null

# Change 2
$find2 = "only.  Four OTO-"
$repl2 = "only.  A single 127mm and four OTO-"
Replace-Text $find2 $repl2 | Out-Null

# Change 3
$find3 = "1993) are were designed to be improved versions"
$repl3 = "1993) were designed as improved versions"
Replace-Text $find3 $repl3 | Out-Null

# Change 4
$find4 = "They have essentially the same armament as the Audace class but with more modern sensors and control systems. One of the 76mm gun has been removed, leaving three. These are well rounded world class ships"
$repl4 = "They have essentially the same armament but with more modern sensors and control systems than their predecessors but with only three 76mm guns vice four. These are well rounded world class ships"
Replace-Text $find4 $repl4 | Out-Null

# Change 5: relocate the _GoBack bookmark
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$markerRange = $d.Content
$markerFound = $markerRange.Find.Execute("vice four", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($markerFound) {
    $markerRange.Collapse(0)   # wdCollapseEnd
    $d.Bookmarks.Add("_GoBack", $markerRange)
} else {
    Write-Host "WARNING - could not find 'vice four' to place _GoBack bookmark"
}

Write-Host "Done"
